# Invalidate formula caches after render.
#
# The "streams|tabular" placeholder block (rows 18-20) gets expanded by the
# templating step, and row 21 holds summary formulas that aggregate over
# that templated range. Previously only K21 existed (a stray empty styled
# cell); now every column A:K gets its own summary formula so that, once a
# consuming application (re)opens the workbook, those formulas recompute
# against whatever data actually ended up in the templated rows instead of
# showing stale/mismatched cached results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A:D summarise the sample numeric columns, E counts the (mostly blank)
# text-ish column, and F:K (the still-empty templated columns) all reduce to
# #DIV/0! until real data lands there.
$ws.Range("A21").Formula = "=AVERAGE(A18:A20)"
$ws.Range("B21").Formula = "=AVERAGE(B18:B20)"
$ws.Range("C21").Formula = "=AVERAGE(C18:C20)"
$ws.Range("D21").Formula = "=AVERAGE(D18:D20)"
$ws.Range("E21").Formula = "=COUNT(E18:E20)"
$ws.Range("F21").Formula = "=AVERAGE(F18:F20)"
$ws.Range("G21").Formula = "=AVERAGE(G18:G20)"
$ws.Range("H21").Formula = "=AVERAGE(H18:H20)"
$ws.Range("I21").Formula = "=AVERAGE(I18:I20)"
$ws.Range("J21").Formula = "=AVERAGE(J18:J20)"

# K21 already existed (empty, with the bold "K column" style picked up from
# the row below) -- reset it back to the plain/general formatting the rest
# of the new summary row uses before giving it the same kind of formula.
$ws.Range("K21").NumberFormat = "General"
$ws.Range("K21").Formula = "=AVERAGE(K18:K20)"
